$wb = $excel.ActiveWorkbook

# Update the comparison text on Tabelle1!C2 to match the style used elsewhere ("vfffff")
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws1.Range("C2").Value = "vfffff"

# Make Tabelle1 the active sheet/tab and select the edited cell, C2
$ws1.Activate()
$ws1.Range("C2").Select()
